$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 (ID 4): add new review entry (fill order matches shared-string allocation order)
$ws.Range("C15").Value = "Défaut+ Cosmétique"
$ws.Range("D15").Value = "Les retours à la ligne ne sont pas respecter selon PEP8 (les espcae entre les fonctions et methodes )"
$ws.Range("B15").Value = "Partout"
$ws.Rows.Item(15).RowHeight = 18.75

# Row 16 (ID 6 -> 5): add new review entry
$ws.Range("A16").Value = 5
$ws.Range("C16").Value = "Défaut+ Cosmétique"
$ws.Range("D16").Value = "les noms des fonctions  (splacher et splasher)"
$ws.Range("B16").Value = "43+48+104"

# Row 17 (ID 7 -> 6): add new review entry
$ws.Range("A17").Value = 6
$ws.Range("C17").Value = "défaut"
$ws.Range("D17").Value = "Toutes les fonctions sont appellés dans programme-principal"
$ws.Range("B17").Value = 128

# Rows 18 through 36: renumber ID column (A) down by 1
for ($r = 18; $r -le 36; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value()
    $cell.Value = $current - 1
}

# Update sheet view: scroll position (topLeftCell ~ A10) and final selection (D30)
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D30").Select()
